$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36-126 down to 37-127
$ws.Rows(36).Insert()

# Populate the new row 36 with data (copy of fixed columns + new values for date/volume/prices)
$ws.Range("A36").Value = 4
$ws.Range("B36").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C36").Value = "Los Lagos"
$ws.Range("D36").Value = 44414
$ws.Range("E36").Value = 10
$ws.Range("F36").Value = 100112043
$ws.Range("G36").Value = "Pepino ensalada"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 400
$ws.Range("K36").Value = 19000
$ws.Range("L36").Value = 19000
$ws.Range("M36").Value = 19000
$ws.Range("N36").Value = "$/caja 60 unidades"
$ws.Range("O36").Value = "Región de Arica y Parinacota"
$ws.Range("P36").Value = 317
$ws.Range("Q36").Value = 60
$ws.Range("R36").Value = "Hortaliza"

# Apply the same style as the other date cells in column D (style index 2 in original workbook)
$ws.Range("D36").Style = $ws.Range("D37").Style
$ws.Range("D36").NumberFormat = $ws.Range("D37").NumberFormat
